# Modify TIM to address DH Study
# - Residential and Services sector have enhanced DH options
# - Four DH demand groups: Low, Medium, High and City
# - New DH Supply options: Surplus heat, heat pumps, solar, & storage.
# - DH cost distribution from Irish Heat Atlas
#
# This adds the 26 Irish county/sub-national "IE-xx" breakdown columns to the
# four RSD* timeslice-data tables (rsdsol, rsd_sh, rsd_rtft, rsd_oe_dem),
# mirroring the column layout already used by tra_dem / other *_dem tables.
# Each new column simply repeats the value already present in the existing
# "National" column of that table's single data row.

$wb = $excel.ActiveWorkbook

# Ordered list of the 26 county columns to append to each table, in the
# exact order used throughout the workbook (see tra_dem / srv_cs_dem tables).
$countyCols = @(
    "IE-CW",
    "IE-KK",
    "IE-LS",
    "IE-LD",
    "IE-LH",
    "IE-OY",
    "IE-WH",
    "IE-WX",
    "IE-CE",
    "IE-KY",
    "IE-TA",
    "IE-LM",
    "IE-MO",
    "IE-RN",
    "IE-SO",
    "IE-CN",
    "IE-DL",
    "IE-MN",
    "IE-D",
    "IE-KE",
    "IE-MH",
    "IE-WW",
    "IE-CO",
    "IE-LK",
    "IE-WD",
    "IE-G"
)

# Sheets whose single ListObject needs the county columns appended. The data
# value for the new columns is copied from the last existing column
# ("National") of that table.
$targetSheets = @("RSDSOL", "RSD_SH", "RSD_RTFT", "RSD_OE_DEM")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lo = $ws.ListObjects.Item(1)

    $lastCol = $lo.ListColumns.Item($lo.ListColumns.Count)
    $nationalValue = $lastCol.DataBodyRange.Value()

    foreach ($colName in $countyCols) {
        $col = $lo.ListColumns.Add()
        $col.Range.Cells.Item(1).Value = $colName
        $col.DataBodyRange.Value = $nationalValue
    }
}

# --- View-state touch ups (cosmetic, mirrors the author's saved selection) ---

$wsSol = $wb.Worksheets.Item("RSDSOL")
$wsSol.Range("AI7").Select()

$wsSh = $wb.Worksheets.Item("RSD_SH")
$wsSh.Range("H2:AG2").Select()

$wsRtft = $wb.Worksheets.Item("RSD_RTFT")
$wsRtft.Range("I2:AH2").Select()

$wsOe = $wb.Worksheets.Item("RSD_OE_DEM")
$wsOe.Range("AF5").Select()
# RSD_OE_DEM is the tab the workbook was left on when last saved.
$wsOe.Activate()

$wsTra = $wb.Worksheets.Item("TRA_DEM")
$wsTra.Range("H2:AG2").Select()
